# Add a new "time" variable row below the existing test-case rows.
# Mirrors the formatting of the row above (F10) and fills in F11 with 5,
# then moves the active selection down to F12 (as Excel does after Enter).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F10").Copy()
$ws.Range("F11").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F11").Value = 5

$ws.Range("F12").Select()
